$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Delete the paragraph "Administrator can create an administrator
#    account X" entirely (removed in the target revision).
# ---------------------------------------------------------------------
$deleteIndex = 16
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("Administrator can create an administrator account")) {
        $deleteIndex = $i
        break
    }
}
$d.Paragraphs.Item($deleteIndex).Range.Delete()

# ---------------------------------------------------------------------
# Helper: change the trailing " X" mark of a paragraph into " V",
# keeping the leading space and the mark itself as two separate runs
# (matching how the change reads in the tracked-changes diff: the
# space run survives untouched, the mark character is replaced).
# ---------------------------------------------------------------------
function Replace-MarkWithV($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pEnd = $p.Range.End
    $rMark = $d.Range($pEnd - 2, $pEnd - 1)
    # Toggling a character formatting property forces the replaced text
    # into its own run instead of being re-merged into the run that
    # precedes it.
    $rMark.Font.Bold = $true
    $rMark.Text = "V"
    $rMark2 = $d.Range($pEnd - 2, $pEnd - 1)
    $rMark2.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Helper: move the hidden "_GoBack" bookmark so that it (collapsed)
# sits right after the last character of the given paragraph, just
# before the paragraph mark.
# ---------------------------------------------------------------------
function Move-GoBackBookmark($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pEnd = $p.Range.End
    # Insert a throw-away placeholder character right before the
    # paragraph mark so we get a genuine (non paragraph-mark-adjacent)
    # insertion point to collapse the bookmark onto, then remove it.
    $insertionPoint = $d.Range($pEnd - 1, $pEnd - 1)
    $insertionPoint.InsertAfter("Z")
    $pEnd2 = $d.Paragraphs.Item($paraIndex).Range.End
    $collapsed = $d.Range($pEnd2 - 2, $pEnd2 - 2)
    $d.Bookmarks.Add("_GoBack", $collapsed)
    $placeholder = $d.Range($pEnd2 - 2, $pEnd2 - 1)
    $placeholder.Delete()
}

# ---------------------------------------------------------------------
# 2. Find every remaining paragraph whose text ends in " X" and turn
#    the mark into " V".
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
    if ($txt.EndsWith(" X")) {
        Replace-MarkWithV($i)
    }
}

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark onto the "Member can modify an
#    income/expenditure" paragraph (it previously sat on "Administrator
#    can see every member's budget").
# ---------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("Member can modify an income/expenditure")) {
        $targetIndex = $i
        break
    }
}
Move-GoBackBookmark($targetIndex)

Write-Output "done"
